$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & volume/1h changes).
# Price/volume cells are stored as plain text in the sheet (e.g. "519.39",
# "  +3.30%  "); values that Excel would otherwise auto-detect as numbers
# are forced to text via NumberFormat "@" and then restored to the default
# "Normal" style so no stray number-format style lingers on the cell.

$ws.Range("D2").Value = "57.264.06"
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").Value = "2.327.76"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").Value = "2.354.00"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").Value = "2.752.09"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "57.075.78"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "2.349.75"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("E22").Value = "  +6.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.80%  "
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  +9.31%  "
$ws.Range("D30").Value = "0.0₃0739"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.45%  "
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.968"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("E40").Value = "  +7.02%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.74%  "
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.563"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("E50").Value = "  +4.27%  "
$ws.Range("E51").Value = "  +2.97%  "
